# "Fixed Commerce Unit Tests, started building Membership Provider replacement."
#
# On the "Generic Backlog" sheet, the backlog item in row 13
# ("Fofson (773) 603-9095, Clarence" / "TODO") is removed entirely.
# Deleting the whole row shifts every row below it up by one
# (old rows 14-17 become the new rows 13-16) and updates the
# active selection to the block that was just revealed/selected
# by the user after the delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Delete()

$ws.Range("A14:B16").Select()
